$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: shift week number, and update MARTES/JUEVES values
$ws.Range("A2").Value = 33
$ws.Range("C2").Value = "SI"
$ws.Range("E2").Value = "SI"

# Row 3: shift week number, swap LUNES/MARTES/JUEVES values, update total
$ws.Range("A3").Value = 34
$ws.Range("B3").Value = "SI"
$ws.Range("C3").Value = "NO"
$ws.Range("E3").Value = "NO"
$ws.Range("G3").Value = 80

# Row 4: shift week number, update LUNES value and total
$ws.Range("A4").Value = 35
$ws.Range("B4").Value = "NO"
$ws.Range("G4").Value = 0
